# Apply the 2025-09-21 01:21:40 JST scrape refresh to the "ランサーズ" sheet.
#
# Every row (2-12) gets its "取得日時" (A column) timestamp bumped to the
# new scrape time. In addition, the scrape that ran at this new timestamp
# returned the two "システム開発" listings (5397121 / 5397117) in the
# opposite order from before, so rows 4 and 5 swap their title (B) and
# URL (F) text while every other column stays put.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-21 01:21:40"

# 1. Refresh the "取得日時" timestamp for every data row (2 through 12).
for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}

# 2. Rows 4 and 5 swap their title (B) and URL (F) values.
$ws.Cells.Item(4, 2).Value = "システム開発において活躍できる案件紹介"
$ws.Cells.Item(4, 6).Value = "https://www.lancers.jp/work/detail/5397117"

$ws.Cells.Item(5, 2).Value = "システム開発の案件紹介"
$ws.Cells.Item(5, 6).Value = "https://www.lancers.jp/work/detail/5397121"
